$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.943.89"
$ws.Range("E2").Value = "  +3.78%  "

$ws.Range("D3").Value = "2.654.80"
$ws.Range("E3").Value = "  +6.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.08"
$ws.Range("E5").Value = "  +2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.82"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("E7").Value = "  +0.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.72"
$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.65"
$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.31"
$ws.Range("E14").Value = "  +2.84%  "

$ws.Range("D15").Value = "3.066.82"
$ws.Range("E15").Value = "  +5.85%  "

$ws.Range("D16").Value = "2.614.95"
$ws.Range("E16").Value = "  +4.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.882"
$ws.Range("E17").Value = "  +5.66%  "

$ws.Range("D18").Value = "49.877.16"
$ws.Range("E18").Value = "  +3.93%  "

$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.07"
$ws.Range("E19").Value = "  +10.34%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.33"
$ws.Range("E20").Value = "  +2.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "282.17"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.21"
$ws.Range("E24").Value = "  +2.40%  "

$ws.Range("E25").Value = "  +2.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.92"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("E28").Value = "  +6.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.74"
$ws.Range("E29").Value = "  +3.89%  "

$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.78"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.58"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("E34").Value = "  +2.55%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0799"
$ws.Range("E36").Value = "  +2.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.08"
$ws.Range("E37").Value = "  +6.91%  "

$ws.Range("E38").Value = "  +2.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +8.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.21"
$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.40"
$ws.Range("E42").Value = "  +5.41%  "

$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("E45").Value = "  +7.30%  "

$ws.Range("D46").Value = "2.070.10"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("E47").Value = "  +14.15%  "

$ws.Range("E48").Value = "  +8.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.41"
$ws.Range("E50").Value = "  +4.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.73"
$ws.Range("E51").Value = "  +1.81%  "
